$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-07"

# Update the header label in column B (shared string) to match the new date
$ws.Range("B1").Value = "October 2021 (through October 07)"

# Row 2 - Garfield Park
$ws.Range("B2").Value = 4
$ws.Range("L2").Value = 6

# Row 4 - North Lawndale
$ws.Range("B4").Value = 3

# Row 6 - Auburn Gresham: new data point
$ws.Range("AF6").Value = 1

# Row 9 - Grand Crossing
$ws.Range("B9").Value = 3

# Row 10 - Roseland: new data point
$ws.Range("AZ10").Value = 1

# Row 38 - Bucktown: new data point
$ws.Range("L38").Value = 1

# Row 48 - Washington Heights
$ws.Range("B48").Value = 2

# Row 65 - Calumet Heights
$ws.Range("B65").Value = 2

# Row 94 - Riverdale: new data point
$ws.Range("B94").Value = 1
